$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.Contains(",") -and $text.Contains("System")) {
            $parts = $text.Split(",")
            for ($i = 0; $i -lt $parts.Length; $i++) {
                $parts[$i] = $parts[$i].Trim()
            }
            $reversed = @()
            for ($i = $parts.Length - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
